$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.921.89"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.642.08"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.71"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5049"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2575"
$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06398"
$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.62"
$ws.Range("E10").Value = "  +0.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07792"
$ws.Range("E11").Value = "  +1.13%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.668.15"
$ws.Range("E12").Value = "  +1.65%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.285"
$ws.Range("E13").Value = "  +1.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5438"
$ws.Range("E14").Value = "  -0.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7866"
$ws.Range("E15").Value = "  -0.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.88"
$ws.Range("E16").Value = "  +2.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.967.16"
$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.005"
$ws.Range("E18").Value = "  -0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.69"
$ws.Range("E19").Value = "  -2.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.398"
$ws.Range("E20").Value = "  +2.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.966"
$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.005"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.007"
$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.868"
$ws.Range("E24").Value = "  -3.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.51"
$ws.Range("E25").Value = "  -0.83%  "

$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.854"
$ws.Range("E27").Value = "  +1.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.74"
$ws.Range("E28").Value = "  -0.14%  "

$ws.Range("E29").Value = "  +0.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04983"
$ws.Range("E30").Value = "  -1.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.265"
$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.197"
$ws.Range("E32").Value = "  +0.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.533"
$ws.Range("E33").Value = "  -0.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.374"
$ws.Range("E34").Value = "  +1.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8942"
$ws.Range("E35").Value = "  +0.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.614"
$ws.Range("E36").Value = "  -0.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.143.51"
$ws.Range("E37").Value = "  -0.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5546"
$ws.Range("E38").Value = "  -1.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01564"
$ws.Range("E39").Value = "  -0.45%  "

$ws.Range("E40").Value = "  -0.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.676"
$ws.Range("E41").Value = "  +0.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8237"
$ws.Range("E42").Value = "  +1.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.81"
$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("E44").Value = "  +8.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.780.33"
$ws.Range("E45").Value = "  +0.22%  "

$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("E47").Value = "  +0.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  -0.31%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05063"
$ws.Range("E49").Value = "  +0.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.008"
$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09511"
$ws.Range("E51").Value = "  +2.25%  "
